$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 464, pushing existing rows 464-494 down to 465-495.
$ws.Rows.Item(464).Insert()

# Populate the new row 464 with the new record (matches the formatting/
# constants shared by the rest of this "Pina" / Macroferia Regional de
# Talca block).
$ws.Range("A464").Value = 5
$ws.Range("B464").Value = "Macroferia Regional de Talca"
$ws.Range("C464").Value = "Maule"
$ws.Range("D464").Value = 45265
$ws.Range("E464").Value = 7
$ws.Range("F464").Value = "Fruta"
$ws.Range("G464").Value = 100108
$ws.Range("H464").Value = "Tropicales y subtropicales"
$ws.Range("I464").Value = 100108005
$ws.Range("J464").Value = "Piña"
$ws.Range("K464").Value = "Caramelo"
$ws.Range("L464").Value = "Segunda"
$ws.Range("M464").Value = 100
$ws.Range("N464").Value = 22000
$ws.Range("O464").Value = 22000
$ws.Range("P464").Value = 22000
$ws.Range("Q464").Value = "$/caja 14 unidades"
$ws.Range("R464").Value = "Ecuador"
$ws.Range("S464").Value = 1571
$ws.Range("T464").Value = 14
